# Auto-generated: apply scheduled-runner market data refresh to Ultima_Profits (Sheets ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(76, 8).Value = 6467.0645
$ws.Cells.Item(76, 9).Value = 5616.6665
$ws.Cells.Item(76, 10).Value = 7004.1577
$ws.Cells.Item(76, 11).Value = 5616.6665
$ws.Cells.Item(76, 12).Value = 7004.1577
$ws.Cells.Item(76, 13).Value = -5301.6665
$ws.Cells.Item(76, 14).Value = -7634.1577
$ws.Cells.Item(79, 8).Value = 6467.0645
$ws.Cells.Item(79, 9).Value = 5616.6665
$ws.Cells.Item(79, 10).Value = 7004.1577
$ws.Cells.Item(79, 11).Value = 5616.6665
$ws.Cells.Item(79, 12).Value = 7004.1577
$ws.Cells.Item(79, 13).Value = -4524.6665
$ws.Cells.Item(79, 14).Value = -9188.1577
$ws.Cells.Item(88, 8).Value = 2570.7144
$ws.Cells.Item(88, 9).Value = 3750
$ws.Cells.Item(88, 10).Value = 998.3333
$ws.Cells.Item(88, 11).Value = 3750
$ws.Cells.Item(88, 12).Value = 998.3333
$ws.Cells.Item(88, 13).Value = -3344
$ws.Cells.Item(88, 14).Value = -1810.3333
$ws.Cells.Item(91, 8).Value = 2570.7144
$ws.Cells.Item(91, 9).Value = 3750
$ws.Cells.Item(91, 10).Value = 998.3333
$ws.Cells.Item(91, 11).Value = 3750
$ws.Cells.Item(91, 12).Value = 998.3333
$ws.Cells.Item(91, 13).Value = -2346
$ws.Cells.Item(91, 14).Value = -3806.3333
$ws.Cells.Item(103, 8).Value = 4623765
$ws.Cells.Item(103, 9).Value = 30050002
$ws.Cells.Item(103, 10).Value = 812.7273
$ws.Cells.Item(103, 11).Value = 90150006
$ws.Cells.Item(103, 12).Value = 2438.1819
$ws.Cells.Item(103, 13).Value = -90149420
$ws.Cells.Item(103, 14).Value = -3610.1819
$ws.Cells.Item(137, 8).Value = 6061950.5
$ws.Cells.Item(137, 9).Value = 888.1
$ws.Cells.Item(137, 10).Value = 15386662
$ws.Cells.Item(137, 11).Value = 2664.3
$ws.Cells.Item(137, 12).Value = 46159986
$ws.Cells.Item(137, 13).Value = -114.3000000000002
$ws.Cells.Item(137, 14).Value = -46165086
$ws.Cells.Item(140, 8).Value = 100133.336
$ws.Cells.Item(140, 10).Value = 100133.336
$ws.Cells.Item(140, 12).Value = 100133.336
$ws.Cells.Item(140, 14).Value = -110493.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9706.6
$ws.Cells.Item(32, 9).Value = 9668.724
$ws.Cells.Item(32, 10).Value = 10300
$ws.Cells.Item(32, 11).Value = 9668.724
$ws.Cells.Item(32, 12).Value = 10300
$ws.Cells.Item(32, 13).Value = -9381.724
$ws.Cells.Item(32, 14).Value = -10874
$ws.Cells.Item(97, 8).Value = 9001.916999999999
$ws.Cells.Item(97, 9).Value = 10621.2
$ws.Cells.Item(97, 10).Value = 905.5
$ws.Cells.Item(97, 11).Value = 10621.2
$ws.Cells.Item(97, 12).Value = 905.5
$ws.Cells.Item(97, 13).Value = -10125.2
$ws.Cells.Item(97, 14).Value = -1897.5
$ws.Cells.Item(139, 8).Value = 42843.438
$ws.Cells.Item(139, 10).Value = 42843.438
$ws.Cells.Item(139, 12).Value = 42843.438
$ws.Cells.Item(139, 14).Value = -53123.438

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 20835152
$ws.Cells.Item(86, 9).Value = 1729.2
$ws.Cells.Item(86, 10).Value = 55557524
$ws.Cells.Item(86, 11).Value = 1729.2
$ws.Cells.Item(86, 12).Value = 55557524
$ws.Cells.Item(86, 13).Value = -606.2
$ws.Cells.Item(86, 14).Value = -55559770
$ws.Cells.Item(89, 8).Value = 20835152
$ws.Cells.Item(89, 9).Value = 1729.2
$ws.Cells.Item(89, 10).Value = 55557524
$ws.Cells.Item(89, 11).Value = 8646
$ws.Cells.Item(89, 12).Value = 277787620
$ws.Cells.Item(89, 13).Value = -3030
$ws.Cells.Item(89, 14).Value = -277798852
$ws.Cells.Item(94, 8).Value = 1305.6666
$ws.Cells.Item(94, 9).Value = 1092.4615
$ws.Cells.Item(94, 10).Value = 1860
$ws.Cells.Item(94, 11).Value = 1092.4615
$ws.Cells.Item(94, 12).Value = 1860
$ws.Cells.Item(94, 13).Value = -641.4614999999999
$ws.Cells.Item(94, 14).Value = -2762

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 886.2
$ws.Cells.Item(107, 9).Value = 739.3
$ws.Cells.Item(107, 10).Value = 1180
$ws.Cells.Item(107, 11).Value = 739.3
$ws.Cells.Item(107, 12).Value = 1180
$ws.Cells.Item(107, 13).Value = 1180.7
$ws.Cells.Item(107, 14).Value = -5020

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(43, 8).Value = 30000
$ws.Cells.Item(43, 10).Value = 30000
$ws.Cells.Item(43, 12).Value = 90000
$ws.Cells.Item(43, 14).Value = -90228
$ws.Cells.Item(87, 8).Value = 15480
$ws.Cells.Item(87, 9).Value = 7950
$ws.Cells.Item(87, 10).Value = 17362.5
$ws.Cells.Item(87, 11).Value = 23850
$ws.Cells.Item(87, 12).Value = 52087.5
$ws.Cells.Item(87, 13).Value = -22602
$ws.Cells.Item(87, 14).Value = -54583.5
$ws.Cells.Item(90, 8).Value = 15480
$ws.Cells.Item(90, 9).Value = 7950
$ws.Cells.Item(90, 10).Value = 17362.5
$ws.Cells.Item(90, 11).Value = 71550
$ws.Cells.Item(90, 12).Value = 156262.5
$ws.Cells.Item(90, 13).Value = -65310
$ws.Cells.Item(90, 14).Value = -168742.5
$ws.Cells.Item(113, 8).Value = 2219.0908
$ws.Cells.Item(113, 9).Value = 402
$ws.Cells.Item(113, 10).Value = 2900.5
$ws.Cells.Item(113, 11).Value = 1206
$ws.Cells.Item(113, 12).Value = 8701.5
$ws.Cells.Item(113, 13).Value = 964
$ws.Cells.Item(113, 14).Value = -13041.5
$ws.Cells.Item(137, 8).Value = 5054295.5
$ws.Cells.Item(137, 9).Value = 9805421
$ws.Cells.Item(137, 10).Value = 6224.25
$ws.Cells.Item(137, 11).Value = 29416263
$ws.Cells.Item(137, 12).Value = 18672.75
$ws.Cells.Item(137, 13).Value = -29411163
$ws.Cells.Item(137, 14).Value = -28872.75
$ws.Cells.Item(140, 8).Value = 2992.7144
$ws.Cells.Item(140, 9).Value = 2992.7144
$ws.Cells.Item(140, 11).Value = 8978.143199999999
$ws.Cells.Item(140, 13).Value = -3798.143199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(126, 8).Value = 3776
$ws.Cells.Item(126, 9).Value = 2535.8462
$ws.Cells.Item(126, 11).Value = 7607.5386
$ws.Cells.Item(126, 13).Value = -5137.5386
$ws.Cells.Item(138, 8).Value = 57570.57
$ws.Cells.Item(138, 10).Value = 57570.57
$ws.Cells.Item(138, 12).Value = 57570.57
$ws.Cells.Item(138, 14).Value = -67850.57000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(82, 8).Value = 1984.7368
$ws.Cells.Item(82, 9).Value = 1793.9166
$ws.Cells.Item(82, 10).Value = 2311.8572
$ws.Cells.Item(82, 11).Value = 1793.9166
$ws.Cells.Item(82, 12).Value = 2311.8572
$ws.Cells.Item(82, 13).Value = -1432.9166
$ws.Cells.Item(82, 14).Value = -3033.8572
$ws.Cells.Item(85, 8).Value = 1984.7368
$ws.Cells.Item(85, 9).Value = 1793.9166
$ws.Cells.Item(85, 10).Value = 2311.8572
$ws.Cells.Item(85, 11).Value = 1793.9166
$ws.Cells.Item(85, 12).Value = 2311.8572
$ws.Cells.Item(85, 13).Value = -545.9166
$ws.Cells.Item(85, 14).Value = -4807.8572
$ws.Cells.Item(93, 8).Value = 2460
$ws.Cells.Item(93, 9).Value = 2933.3333
$ws.Cells.Item(93, 10).Value = 1750
$ws.Cells.Item(93, 11).Value = 2933.3333
$ws.Cells.Item(93, 12).Value = 1750
$ws.Cells.Item(93, 13).Value = -1685.3333
$ws.Cells.Item(93, 14).Value = -4246
$ws.Cells.Item(139, 8).Value = 41686.727
$ws.Cells.Item(139, 10).Value = 41790.4
$ws.Cells.Item(139, 12).Value = 41790.4
$ws.Cells.Item(139, 14).Value = -52070.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(51, 8).Value = 11035.6
$ws.Cells.Item(51, 10).Value = 11035.6
$ws.Cells.Item(51, 12).Value = 11035.6
$ws.Cells.Item(51, 14).Value = -12055.6
$ws.Cells.Item(52, 8).Value = 7209.4
$ws.Cells.Item(52, 9).Value = 2500
$ws.Cells.Item(52, 10).Value = 10349
$ws.Cells.Item(52, 11).Value = 2500
$ws.Cells.Item(52, 12).Value = 10349
$ws.Cells.Item(52, 13).Value = -2274
$ws.Cells.Item(52, 14).Value = -10801
$ws.Cells.Item(81, 8).Value = 805
$ws.Cells.Item(81, 9).Value = 989.5
$ws.Cells.Item(81, 10).Value = 620.5
$ws.Cells.Item(81, 11).Value = 1979
$ws.Cells.Item(81, 12).Value = 1241
$ws.Cells.Item(81, 13).Value = -918
$ws.Cells.Item(81, 14).Value = -3363
$ws.Cells.Item(84, 8).Value = 805
$ws.Cells.Item(84, 9).Value = 989.5
$ws.Cells.Item(84, 10).Value = 620.5
$ws.Cells.Item(84, 11).Value = 9895
$ws.Cells.Item(84, 12).Value = 6205
$ws.Cells.Item(84, 13).Value = -4591
$ws.Cells.Item(84, 14).Value = -16813
$ws.Cells.Item(138, 8).Value = 64923.332
$ws.Cells.Item(138, 9).Value = 0
$ws.Cells.Item(138, 10).Value = 64923.332
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 64923.332
$ws.Cells.Item(138, 14).Value = -75203.33199999999
$ws.Cells.Item(138, 13).ClearContents()
